$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $new, $wdReplaceAll)
}

Replace-Text "Ativação: 01/01/2019" "Ativação: 01/01/2022"

Replace-Text "Apresentar aos alunos de Engenharia Bioquímica os conceitos básicos de matemática financeira e Engenharia Econômica." "Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia"

Replace-Text "Introdução à Engenharia Econômica; fluxo de caixa; métodos de decisão Depreciação; Excel aplicado à Engenharia Econômica." "Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa; depreciação; comparação entre alternativa de investimentos; ponto de equilíbrio; planilha eletrônica aplicado à Engenharia Econômica."

Replace-Text "1.Introdução à Engenharia Econômica: engenharia econômica como uma ferramenta de análise de tomada de decisão; variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência.2.Fluxo de Caixa: estimativa e representação gráfica.3.Métodos de decisão: taxa mínima de atratividade, critérios econômicos de decisão (método do valor presente, método do valor anual, método da taxa interna de retorno); renovação e substituição de equipamentos. 4.Depreciação: métodos de depreciação.5.Excel aplicado à Engenharia Econômica: Utilização de planilhas eletrônicas." "1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto (engenharia econômica como uma ferramenta de análise de tomada de decisão); 2.ESTIMATIVA DE CUTOS: estimativa de investimento de capital (tipos de estimativas de custo de capital; índice de custos mais comuns; métodos para estimar investimento de capital); estimativa de custo de equipamentos; estimativa de custo de produção;3.JUROS: variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência. 4.FLUXO DE CAIXA: diagrama de fluxo de caixa; equivalência de fluxo de caixa; séries uniforme e em gradientes; elaboração de fluxo de caixa5.Depreciação: métodos de depreciação;6.COMPARAÇÃO ENTRE ALTERNATIVAS DE INVESTIMENTO: critério de Rentabilidade – método do Valor Anual Uniforme Equivalente (VAU); método do Valor Presente (VP); método da Taxa interna de Retorno (TIR); método da Taxa Interna de Retorno Modificada (TIRM); renovação e substituição de equipamentos; período de retorno;7.Ponto de Equilíbrio de um projeto (break-even-point);8.Planilhas eletrônicas aplicado à análises de Engenharia Econômica."

Replace-Text "A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição da notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n." "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Econômica na área de Engenharia Bioquímica, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão."

Replace-Text "A nota Final (NF) corresponderá à Média Aritmética das Avaliações" "A nota do projeto (NP) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, avaliação dos Pares, apresentação de trabalhos, dentre outros;A nota dos exercícios (NE) será a média aritmética dos exercícios individuais;NF = (0,8NP + 0,2NE)"

Replace-Text "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2" "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2n"

Replace-Text "1.Fraser, Niall M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013.2.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009.3.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007.4.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005." "1.Blank L. e Tarquin A. Engenharia Econômica. 6ed. Mc GrawHill, 2007.2.Fraser N. M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013. 3.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009. 4.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007. 5.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005 6.Peter, MS. And Timmerhaus, K.D. Plant Design and economic for chemical engineers"
